$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: status changes from "Chưa Hoàn Thành" to "Hoàn thành"
$ws.Range("F4").Value = "Hoàn thành"

# Rows 5-7: drop the stray params[] ("#") value in column D and add a
# status of "Hoàn thành" in column F
$ws.Range("D5").ClearContents()
$ws.Range("F5").Value = "Hoàn thành"

$ws.Range("D6").ClearContents()
$ws.Range("F6").Value = "Hoàn thành"

$ws.Range("D7").ClearContents()
$ws.Range("F7").Value = "Hoàn thành"

# New row 8: Login route
$ws.Range("A8").Value = "Nhom420/Login"
$ws.Range("B8").Value = "Login"
$ws.Range("C8").Value = "showDefault"
$ws.Range("E8").Value = "Hiện màn hình login"
$ws.Range("F8").Value = "Chưa"

# New row 9: Register route
$ws.Range("A9").Value = "Nhom420/Register"
$ws.Range("B9").Value = "Register"
$ws.Range("C9").Value = "showDefault"
$ws.Range("E9").Value = "Hiện màn hình đăng kí"
$ws.Range("F9").Value = "Chưa"

# Update the recorded selection to match the saved view state
[void]$ws.Range("E14").Select()
